$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.279.06'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '2.515.86'
$ws.Range("E3").Value = '  +2.92%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''542.39'
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("D6").Value = '''144.66'
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").Value = '''0.996'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '''0.575'
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("D9").Value = '2.545.78'
$ws.Range("E9").Value = '  +3.58%  '
$ws.Range("E10").Value = '  +1.86%  '
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("E12").Value = '  +4.62%  '
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D14").Value = '2.963.07'
$ws.Range("D15").Value = '''23.81'
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").Value = '59.242.49'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("E17").Value = '  +2.27%  '
$ws.Range("D18").Value = '2.545.56'
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("D19").Value = '''11.26'
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").Value = '''4.30'
$ws.Range("E20").Value = '  -1.17%  '
$ws.Range("D21").Value = '''326.42'
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("E22").Value = '  +3.09%  '
$ws.Range("D23").Value = '''5.85'
$ws.Range("E23").Value = '  +2.77%  '
$ws.Range("D24").Value = '''62.06'
$ws.Range("E24").Value = '  +2.30%  '
$ws.Range("D25").Value = '''0.438'
$ws.Range("E25").Value = '  -2.77%  '
$ws.Range("E26").Value = '  +2.63%  '
$ws.Range("D27").Value = '''0.992'
$ws.Range("E27").Value = '  +1.61%  '
$ws.Range("D28").Value = '''8.05'
$ws.Range("E28").Value = '  +5.06%  '
$ws.Range("D29").Value = '''6.87'
$ws.Range("E29").Value = '  +3.62%  '
$ws.Range("D30").Value = '0.0₃0784'
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("D31").Value = '''1.83'
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("E33").Value = '  +8.85%  '
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = '''156.80'
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").Value = '''18.68'
$ws.Range("E36").Value = '  +1.56%  '
$ws.Range("D37").Value = '''4.38'
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("D38").Value = '''1.61'
$ws.Range("E38").Value = '  -5.06%  '
$ws.Range("D39").Value = '''5.64'
$ws.Range("E39").Value = '  -3.05%  '
$ws.Range("D40").Value = '''36.94'
$ws.Range("E40").Value = '  +2.15%  '
$ws.Range("D41").Value = '''299.06'
$ws.Range("E41").Value = '  -4.80%  '
$ws.Range("D42").Value = '''3.71'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").Value = '''0.829'
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").Value = '''0.604'
$ws.Range("E45").Value = '  +4.09%  '
$ws.Range("D46").Value = '''10.80'
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("D47").Value = '''0.0936'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = '''18.85'
$ws.Range("E48").Value = '  +2.82%  '
$ws.Range("D49").Value = '''123.88'
$ws.Range("E49").Value = '  +1.66%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").Value = '''0.0515'
$ws.Range("E51").Value = '  -1.65%  '
